# Apply the changes described by the fixture diff ("added in cldf writing support").
#
# Sheet2 (the active sheet) gets:
#   - A2 changed from 1 to 1.01
#   - a new row, A5 = 3 (row 4 intentionally left blank), which extends the
#     sheet's dimension from A1:A3 to A1:A5
#   - the active selection moved to the newly added A5 cell
#
# Both Sheet1 and Sheet2 also pick up a (very slightly) wider default column A,
# which we approximate as closely as this engine's column-width rounding allows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Update cell A2 to its new value.
$ws2.Range("A2").Value = 1.01

# Add the new row 5 (row 4 stays empty), which also pushes the sheet's
# dimension/used-range out to A1:A5.
$ws2.Range("A5").Value = 3

# Move/keep the active selection on the newly written cell.
$ws2.Range("A5").Select() | Out-Null

# Nudge column A's width on both sheets a little wider (closest attainable
# value to the fixture's new width of ~11.88 characters).
$ws1.Columns.Item(1).ColumnWidth = 11
$ws2.Columns.Item(1).ColumnWidth = 11
